# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (totals) sheet,
#    populated with the per-fund holding detail for that quarter.
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet (date, holding
#    count, holding market value), re-indexing the existing rows.

$wb = $excel.ActiveWorkbook

# Grab stable, name-based handles before doing any structural changes.
# (Worksheets.Item(<number>) in this runtime re-resolves by position, so a
# variable captured that way goes stale once sheets are inserted/reordered;
# name-based lookups stay valid no matter how the tab order changes.)
$styleSrcName = $wb.Worksheets.Item($wb.Worksheets.Count - 1).Name   # "2021-Q4"
$totalSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name     # "总计"

$totalSheet = $wb.Worksheets.Item($totalSheetName)

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" detail sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Re-resolve by name now that the tab order has changed.
$styleSrc = $wb.Worksheets.Item($styleSrcName)
$totalSheet = $wb.Worksheets.Item($totalSheetName)
$q1 = $wb.Worksheets.Item("2022-Q1")

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund holding detail rows: index, code, name, size, stock-position, position-ratio, market-value, rank
$fundRows = @(
    @(0, "009312", "新疆前海联合价值优选混合A",     "10.35", "92.42", "5.03", "0.5206", 6),
    @(1, "004693", "新疆前海联合泳隽灵活配置混合A", "9.08",  "93.74", "5.00", "0.4540", 5),
    @(2, "001146", "中欧瑾源灵活配置混合 - A",       "8.19",  "23.06", "1.93", "0.1581", 3),
    @(3, "007066", "浦银安盛先进制造混合A",          "3.66",  "74.55", "4.15", "0.1519", 6),
    @(4, "007067", "浦银安盛先进制造混合C",          "2.27",  "74.55", "4.15", "0.0942", 6),
    @(5, "009313", "新疆前海联合价值优选混合C",      "1.67",  "92.42", "5.03", "0.0840", 6),
    @(6, "001147", "中欧瑾源灵活配置混合 - C",       "3.82",  "23.06", "1.93", "0.0737", 3),
    @(7, "004734", "中欧瑾灵灵活配置混合A",          "3.55",  "32.29", "1.79", "0.0635", 6),
    @(8, "004735", "中欧瑾灵灵活配置混合C",          "0.33",  "32.29", "1.79", "0.0059", 6),
    @(9, "007042", "新疆前海联合泳隽灵活配置混合C",  "0.00",  "93.74", "5.00", $null,    5)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 1).Value = $row[0]          # A: index (number)
    $q1.Cells.Item($r, 2).Value = "'" + $row[1]    # B: fund code (keep as text, keep leading zeros)
    $q1.Cells.Item($r, 3).Value = $row[2]          # C: fund name (text)
    $q1.Cells.Item($r, 4).Value = "'" + $row[3]    # D: fund size (text)
    $q1.Cells.Item($r, 5).Value = "'" + $row[4]    # E: total stock position (text)
    $q1.Cells.Item($r, 6).Value = "'" + $row[5]    # F: position ratio (text)
    if ($row[6] -eq $null) {
        $q1.Cells.Item($r, 7).Value = 0            # G: holding market value - real 0 for the last row
    } else {
        $q1.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $q1.Cells.Item($r, 8).Value = $row[7]          # H: position rank (number)
    $r = $r + 1
}

# Apply the same header/index styling used by the other quarterly sheets
$styleSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$styleSrc.Range("A2").Copy()
$q1.Range("A2:A11").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. "总计" sheet - prepend the 2022-Q1 summary row
# ---------------------------------------------------------------------------
$totalSheet.Cells.Clear()

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$summaryRows = @(
    @("2022-Q1", 10, 1.61),
    @("2021-Q4", 32, 10.64),
    @("2021-Q3", 22, 7.07),
    @("2021-Q2", 33, 17.69),
    @("2021-Q1", 71, 38.02),
    @("2020-Q4", 73, 39.18)
)

$r = 2
$idx = 0
foreach ($row in $summaryRows) {
    $totalSheet.Cells.Item($r, 1).Value = $idx
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}

$styleSrc.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)

$styleSrc.Range("A2").Copy()
$totalSheet.Range("A2:A7").PasteSpecial(-4122)
